# Auto-generated Excel COM-interop script
# Applies updated market price data cells across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2999.4
$ws.Range("J40").Value = 2999.4
$ws.Range("L40").Value = 2999.4
$ws.Range("N40").Value = -3349.4
$ws.Range("H112").Value = 1560.0769
$ws.Range("I112").Value = 987
$ws.Range("J112").Value = 1814.7778
$ws.Range("K112").Value = 2961
$ws.Range("L112").Value = 5444.3334
$ws.Range("M112").Value = -1853
$ws.Range("N112").Value = -7660.3334
$ws.Range("H137").Value = 2211.4285
$ws.Range("I137").Value = 2046
$ws.Range("J137").Value = 2432
$ws.Range("K137").Value = 6138
$ws.Range("L137").Value = 7296
$ws.Range("M137").Value = -3588
$ws.Range("N137").Value = -12396
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2599.4443
$ws.Range("I45").Value = 1959.2
$ws.Range("K45").Value = 1959.2
$ws.Range("M45").Value = -1582.2
$ws.Range("H61").Value = 3499.25
$ws.Range("I61").Value = 3499.25
$ws.Range("K61").Value = 3499.25
$ws.Range("M61").Value = -3287.25
$ws.Range("H122").Value = 1711499.1
$ws.Range("I122").Value = 2539998.8
$ws.Range("J122").Value = 54500
$ws.Range("K122").Value = 7619996.399999999
$ws.Range("L122").Value = 163500
$ws.Range("M122").Value = -7617546.399999999
$ws.Range("N122").Value = -168400
$ws.Range("H132").Value = 3558.0454
$ws.Range("I132").Value = 3487.1765
$ws.Range("K132").Value = 10461.5295
$ws.Range("M132").Value = -7931.529500000001
$ws.Range("H136").Value = 3499.25
$ws.Range("I136").Value = 3499.25
$ws.Range("K136").Value = 10497.75
$ws.Range("M136").Value = -7947.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4337.75
$ws.Range("I20").Value = 3150.2
$ws.Range("J20").Value = 5186
$ws.Range("K20").Value = 3150.2
$ws.Range("L20").Value = 5186
$ws.Range("M20").Value = -2903.2
$ws.Range("N20").Value = -5680
$ws.Range("H42").Value = 150000
$ws.Range("J42").Value = 150000
$ws.Range("L42").Value = 150000
$ws.Range("N42").Value = -150656
$ws.Range("H99").Value = 2721.8333
$ws.Range("I99").Value = 2721.8333
$ws.Range("K99").Value = 2721.8333
$ws.Range("M99").Value = -1223.8333
$ws.Range("H105").Value = 3090183.2
$ws.Range("I105").Value = 4905003
$ws.Range("K105").Value = 4905003
$ws.Range("M105").Value = -4903256
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71431490
$ws.Range("I16").Value = 125002750
$ws.Range("J16").Value = 3133
$ws.Range("K16").Value = 125002750
$ws.Range("L16").Value = 3133
$ws.Range("M16").Value = -125002463
$ws.Range("N16").Value = -3707
$ws.Range("H58").Value = 3912.625
$ws.Range("I58").Value = 2270.6667
$ws.Range("J58").Value = 4897.8
$ws.Range("K58").Value = 2270.6667
$ws.Range("L58").Value = 4897.8
$ws.Range("M58").Value = -2067.6667
$ws.Range("N58").Value = -5303.8
$ws.Range("H62").Value = 102001
$ws.Range("I62").Value = 3002.5
$ws.Range("J62").Value = 200999.5
$ws.Range("K62").Value = 3002.5
$ws.Range("L62").Value = 200999.5
$ws.Range("M62").Value = -2378.5
$ws.Range("N62").Value = -202247.5
$ws.Range("H65").Value = 102001
$ws.Range("I65").Value = 3002.5
$ws.Range("J65").Value = 200999.5
$ws.Range("K65").Value = 15012.5
$ws.Range("L65").Value = 1004997.5
$ws.Range("M65").Value = -11892.5
$ws.Range("N65").Value = -1011237.5
$ws.Range("H94").Value = 3699.4
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 3874.25
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 3874.25
$ws.Range("M94").Value = -2549
$ws.Range("N94").Value = -4776.25
$ws.Range("H113").Value = 71431490
$ws.Range("I113").Value = 125002750
$ws.Range("J113").Value = 3133
$ws.Range("K113").Value = 125002750
$ws.Range("L113").Value = 3133
$ws.Range("M113").Value = -125000580
$ws.Range("N113").Value = -7473
$ws.Range("H134").Value = 2452.818
$ws.Range("I134").Value = 1711.8572
$ws.Range("J134").Value = 3749.5
$ws.Range("K134").Value = 5135.571599999999
$ws.Range("L134").Value = 11248.5
$ws.Range("M134").Value = -2600.571599999999
$ws.Range("N134").Value = -16318.5
$ws.Range("H136").Value = 3912.625
$ws.Range("I136").Value = 2270.6667
$ws.Range("J136").Value = 4897.8
$ws.Range("K136").Value = 6812.000100000001
$ws.Range("L136").Value = 14693.4
$ws.Range("M136").Value = -4262.000100000001
$ws.Range("N136").Value = -19793.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 688.5
$ws.Range("I5").Value = 899.5
$ws.Range("J5").Value = 477.5
$ws.Range("K5").Value = 2698.5
$ws.Range("L5").Value = 1432.5
$ws.Range("M5").Value = -2586.5
$ws.Range("N5").Value = -1656.5
$ws.Range("H23").Value = 500099.5
$ws.Range("J23").Value = 500099.5
$ws.Range("L23").Value = 1500298.5
$ws.Range("N23").Value = -1500768.5
$ws.Range("H59").Value = 5000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = $null
$ws.Range("N59").Value = -16080
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = $null
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H97").Value = 5124
$ws.Range("J97").Value = 6165.3335
$ws.Range("L97").Value = 18496.0005
$ws.Range("N97").Value = -19488.0005
$ws.Range("H104").Value = 8333
$ws.Range("J104").Value = 9999.5
$ws.Range("L104").Value = 29998.5
$ws.Range("N104").Value = -35240.5
$ws.Range("H135").Value = 688.5
$ws.Range("I135").Value = 899.5
$ws.Range("J135").Value = 477.5
$ws.Range("K135").Value = 8095.5
$ws.Range("L135").Value = 4297.5
$ws.Range("M135").Value = -5560.5
$ws.Range("N135").Value = -9367.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H129").Value = 29999.5
$ws.Range("J129").Value = 29999.5
$ws.Range("L129").Value = 29999.5
$ws.Range("N129").Value = -39999.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1125.25
$ws.Range("I22").Value = 1167
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1167
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -872
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1125.25
$ws.Range("I27").Value = 1167
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1167
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -1060
$ws.Range("N27").Value = -1214
$ws.Range("H61").Value = 37043036
$ws.Range("I61").Value = 37043036
$ws.Range("K61").Value = 37043036
$ws.Range("M61").Value = -37042834
$ws.Range("H113").Value = 37043036
$ws.Range("I113").Value = 37043036
$ws.Range("K113").Value = 37043036
$ws.Range("M113").Value = -37040866
$ws.Range("H124").Value = 54607
$ws.Range("J124").Value = 54607
$ws.Range("L124").Value = 54607
$ws.Range("N124").Value = -64427
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3707.4
$ws.Range("I62").Value = 2471.6667
$ws.Range("J62").Value = 4402.5
$ws.Range("K62").Value = 2471.6667
$ws.Range("L62").Value = 4402.5
$ws.Range("M62").Value = -1847.6667
$ws.Range("N62").Value = -5650.5
$ws.Range("H65").Value = 3707.4
$ws.Range("I65").Value = 2471.6667
$ws.Range("J65").Value = 4402.5
$ws.Range("K65").Value = 12358.3335
$ws.Range("L65").Value = 22012.5
$ws.Range("M65").Value = -9238.333500000001
$ws.Range("N65").Value = -28252.5
$ws.Range("H81").Value = 8703.637000000001
$ws.Range("I81").Value = 5105.7144
$ws.Range("K81").Value = 10211.4288
$ws.Range("M81").Value = -9150.4288
$ws.Range("H84").Value = 8703.637000000001
$ws.Range("I84").Value = 5105.7144
$ws.Range("K84").Value = 51057.144
$ws.Range("M84").Value = -45753.144
$ws.Range("H132").Value = 2466.9
$ws.Range("I132").Value = 2132.32
$ws.Range("K132").Value = 6396.960000000001
$ws.Range("M132").Value = -3866.960000000001
$ws.Range("H135").Value = 79998.5
$ws.Range("J135").Value = 79998.5
$ws.Range("L135").Value = 79998.5
$ws.Range("N135").Value = -90138.5
